$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("G3").Value = 68
$ws1.Range("F6").Value = 149
$ws1.Range("F7").Value = 156
$ws1.Range("F8").Value = 4594
$ws1.Range("F12").Value = 460
$ws1.Range("F13").Value = 18
$ws1.Range("F15").Value = 1325
$ws1.Range("F16").Value = 2567
$ws1.Range("F18").Value = 78
$ws1.Range("F19").Value = 54
$ws1.Range("F21").Value = 2290
$ws1.Range("F24").Value = 27
$ws1.Range("F25").Value = 160
$ws1.Range("F26").Value = 105
$ws1.Range("F28").Value = 225
$ws1.Range("F29").Value = 36

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G3").Value = 68
$ws4.Range("F6").Value = 149
$ws4.Range("F7").Value = 156
$ws4.Range("F8").Value = 39
$ws4.Range("F9").Value = 4594
$ws4.Range("F13").Value = 460
$ws4.Range("F14").Value = 18
$ws4.Range("F16").Value = 1325
$ws4.Range("F17").Value = 2567
$ws4.Range("F19").Value = 78
$ws4.Range("F20").Value = 54
$ws4.Range("F22").Value = 2290
$ws4.Range("F25").Value = 27
$ws4.Range("F26").Value = 160
$ws4.Range("F27").Value = 105
$ws4.Range("F29").Value = 225
$ws4.Range("F30").Value = 36
